$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 16: re-style the table (a:tableStyleId change)
#    {9D70FA5D-22BF-47B0-9787-609BD44AD113} -> {DF77905B-2622-4569-AB46-0C2B584F6C6A}
# ---------------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
for ($i = 1; $i -le $s16.Shapes.Count; $i++) {
    $shp = $s16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{DF77905B-2622-4569-AB46-0C2B584F6C6A}")
    }
}

# ---------------------------------------------------------------------------
# 2) Theme swap: the deck's two theme parts (theme1.xml = "Office Theme",
#    theme2.xml = "Integral") had their colour schemes exchanged. The live/
#    editable theme reachable from the object model is the one bound to the
#    slide master (theme2.xml, "Integral"); recolour it to the values that
#    "Office Theme" used to carry.
# ---------------------------------------------------------------------------
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x6A5444,  # dk2      (44546A)
    0xE6E6E7,  # lt2      (E7E6E6)
    0xD59B5B,  # accent1  (5B9BD5)
    0x317DED,  # accent2  (ED7D31)
    0xA5A5A5,  # accent3  (A5A5A5)
    0x00C0FF,  # accent4  (FFC000)
    0xC47244,  # accent5  (4472C4)
    0x47AD70,  # accent6  (70AD47)
    0xC16305,  # hlink    (0563C1)
    0x724F95   # folHlink (954F72)
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
